# Regenerate all penyata to follow new data and format
# - Rename "Kali Pertama/Kedua/Ketiga/Keempat" labels to "Semakan Kali ..."
#   (these labels are reused across the three sub-sections: Penandaan Fail,
#   Laporan Atas Talian, and JPPM / JDM / JDRM)
# - Normalise the competition-entry names from ALL CAPS to Title Case
# - Fill in the "Semakan Kali Ketiga" Merit / Demerit figures for the
#   Penandaan Fail section (row 18)
# - Move the "STATEMENT OF HOMEROOM ACCOUNT" title over one column (E4 -> D4)
# - Re-work a handful of merged-cell ranges to match the new layout
# - Centre the printout horizontally, zero out the header/footer margins,
#   and constrain the print area to exactly one page tall

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Label renames (applies to every cell that shares the text) -------
$renameMap = @{
    "Kali Pertama"         = "Semakan Kali Pertama"
    "Kali Kedua"           = "Semakan Kali Kedua"
    "Kali Ketiga"          = "Semakan Kali Ketiga"
    "Kali Keempat"         = "Semakan Kali Keempat"
    "BOUQUET KREATIF"      = "Bouquet Kreatif"
    "TIK TOK RAYA"         = "Tik Tok Raya"
    "RIANG RIA KUIH RAYA"  = "Riang Ria Kuih Raya"
    "CREATIVE COLLAGE"     = "Creative Collage"
}

foreach ($row in 1..47) {
    foreach ($colLetter in @("A","B","C","D","E","F","G")) {
        $cell = $ws.Range("$colLetter$row")
        $val = $cell.Value2
        if ($renameMap.ContainsKey($val)) {
            $cell.Value = $renameMap[$val]
        }
    }
}

# --- 2. Updated Merit / Demerit figures for "Semakan Kali Ketiga" --------
$ws.Range("D18").Value = 9630
$ws.Range("E18").Value = 1340

# --- 3. Move the statement title from E4 to D4 ---------------------------
$titleValue = $ws.Range("E4").Value2
$ws.Range("D4").Value = $titleValue
$ws.Range("E4").Value = $null

# --- 4. Merged-cell layout updates ----------------------------------------
$ws.Range("B15:C15").UnMerge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("D4:G4").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()
$ws.Range("B43:E43").Merge()

# --- 5. Print / page setup -------------------------------------------------
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.HeaderMargin = 0
$ws.PageSetup.FooterMargin = 0
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
